# Lisää API-järjestelmätunnukset osio oikeuksiin
# Inserts a new row into the "Oikeudet" sheet, right after the existing
# "Hallinta / Integraatioloki" row, for the new "API-järjestelmätunnukset"
# view under the "Hallinta" section with "R,W" rights.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

# Row 58 ("Hallinta" / "Integraatioloki") is the last row of the "Hallinta"
# section; push everything from row 58 onward down by one row.
$ws.Rows.Item(58).Insert()

# Match the formatting of the row directly above it (same section / same
# bordered+centered look used throughout the table).
$rng = $ws.Range("A58:W58")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
$rng.Font.Name = "Arial"
$rng.Font.Size = 10
$ws.Range("C58:W58").HorizontalAlignment = -4108
$ws.Rows.Item(58).RowHeight = 15.75

$ws.Cells.Item(58, 1).Value = "Hallinta"
$ws.Cells.Item(58, 2).Value = "API-järjestelmätunnukset"
$ws.Cells.Item(58, 4).Value = "R,W"
